$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 'Isaiah Collier'
$ws.Range("B2").Value = 'PG'
$ws.Range("C2").Value = 'Utah Jazz'
$ws.Range("A3").Value = 'Tyler Herro'
$ws.Range("B3").Value = 'PG,SG'
$ws.Range("C3").Value = 'Miami Heat'
$ws.Range("A4").Value = 'De''Aaron Fox'
$ws.Range("B4").Value = 'PG'
$ws.Range("C4").Value = 'Sacramento Kings'
$ws.Range("A5").Value = 'Mikal Bridges'
$ws.Range("B5").Value = 'SG,SF,PF'
$ws.Range("C5").Value = 'New York Knicks'
$ws.Range("A6").Value = 'Miles Bridges'
$ws.Range("B6").Value = 'SF,PF'
$ws.Range("C6").Value = 'Charlotte Hornets'
$ws.Range("A7").Value = 'DeMar DeRozan'
$ws.Range("B7").Value = 'SF,PF'
$ws.Range("C7").Value = 'Sacramento Kings'
$ws.Range("A8").Value = 'Harrison Barnes'
$ws.Range("B8").Value = 'SF,PF'
$ws.Range("C8").Value = 'San Antonio Spurs'
$ws.Range("A9").Value = 'Shaedon Sharpe'
$ws.Range("B9").Value = 'SG,SF'
$ws.Range("C9").Value = 'Portland Trail Blazers'
$ws.Range("A10").Value = 'Nick Richards'
$ws.Range("B10").Value = 'C'
$ws.Range("C10").Value = 'Phoenix Suns'
$ws.Range("A11").Value = 'T.J. McConnell'
$ws.Range("B11").Value = 'PG'
$ws.Range("C11").Value = 'Indiana Pacers'
$ws.Range("A12").Value = 'Nikola Vucevic'
$ws.Range("B12").Value = 'PF,C'
$ws.Range("C12").Value = 'Chicago Bulls'
$ws.Range("A13").Value = 'Brook Lopez'
$ws.Range("B13").Value = 'C'
$ws.Range("C13").Value = 'Milwaukee Bucks'
$ws.Range("A14").Value = 'Josh Giddey'
$ws.Range("B14").Value = 'PG,SG,SF'
$ws.Range("C14").Value = 'Chicago Bulls'
$ws.Range("A15").Value = 'Scottie Barnes'
$ws.Range("B15").Value = 'PG,SG,SF,PF'
$ws.Range("C15").Value = 'Toronto Raptors'
$ws.Range("A16").Value = 'Evan Mobley'
$ws.Range("B16").Value = 'PF,C'
$ws.Range("C16").Value = 'Cleveland Cavaliers'
$ws.Range("A17").Value = 'Luka Doncic'
$ws.Range("B17").Value = 'PG,SG'
$ws.Range("C17").Value = 'Dallas Mavericks'
$ws.Range("A18").Value = 'Bobby Portis'
$ws.Range("B18").Value = 'PF,C'
$ws.Range("C18").Value = 'Milwaukee Bucks'
$ws.Range("A19").Value = 'Ja Morant'
$ws.Range("B19").Value = 'PG'
$ws.Range("C19").Value = 'Memphis Grizzlies'
Write-Output "done"
